# Re-organizing files & re-doing simulation inputs
# Applies the edits described by the commit to Sheet1 of the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simulation-ended timestamp (B2) is no longer known -> clear it ---
# ("Simulation started:" in B1 is left as-is.)
$ws.Range("B2").Value = ""

# --- Input/Output file names bumped from the "23" run to the "24" run ---
$ws.Range("B7").Value = "input24.csv"
$ws.Range("B8").Value = "output24.csv"

# --- "Conditions varied" values reformatted / redone ---
$ws.Range("C23").Value = "1,15,30,60"
$ws.Range("C24").Value = "1,15,30,60"
$ws.Range("C25").Value = "4, 0.2"

# --- A23/A24 ("initial_perc_SAV_cover"/"initial_perc_FP_cover") and A25
#     ("replicates") gain a thin box border around them ---
foreach ($addr in @("A23", "A24", "A25")) {
    $c = $ws.Range($addr)
    $c.Borders.LineStyle = 1  # xlContinuous
    $c.Borders.Weight = 2     # xlThin
}

# --- Update the active selection to reflect where the editor left off ---
$ws.Range("F46").Select() | Out-Null
